# "new results. 2021/06/03 10:32"
#
# The author reran the experiment and pasted in a refreshed data point for
# expert #29 (B30), which ripples through the AVERAGE() summary in B32.
# They also simply had a different cell selected / scrolled to when they
# last saved (selection moves to C32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New measurement for row 29 ("2 Experts" / B30): 0.4012 -> 0.40600000000000003
$ws.Range("B30").Value = 0.40600000000000003

# B32 = AVERAGE(B2:B31) is a live formula, so it recalculates automatically
# from 0.38927999999999996 to 0.38943999999999995 once B30 changes above -
# no explicit write needed.

# Cursor/selection left on C32 (the averages cell) when the file was saved.
$ws.Range("C32").Select()
